# ---------------------------------------------------------------------------
# Edits comunas_descripcion.xlsx:
#  1. Duplicates "Hoja1" into a second sheet "documentacion_descripciones"
#     that preserves the original, fully-documented table (columns D/E/F
#     with "palabras clave" / "Paginas" / "Refeferencias").
#  2. Renames the original sheet to "descripciones_comunas" and turns it
#     into the simplified/display table: fills in the previously-empty
#     "descripcion" cells, and clears out the no-longer-needed D/E/F
#     content (while keeping the header/row formatting).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Duplicate the sheet first, so the copy keeps ALL the original data
#        (including the D/E/F "documentation" columns) before we touch it. ---
$ws.Copy([System.Type]::Missing, $ws)
$wsDoc = $wb.ActiveSheet
$wsDoc.Name = "documentacion_descripciones"

# Documentation sheet is not the tab that should stay selected, and its
# stored cursor position is C21 (no multi-row header selection).
$wsDoc.Range("C21").Select()

# --- 2. Rename / reselect the original sheet -------------------------------
$ws.Name = "descripciones_comunas"
$ws.Select()
$ws.Range("B1:C1").Select()

# --- 3. Fill in the previously-blank "descripcion" cells on the display
#        sheet with the new shorthand descriptions. Written in the exact
#        order the new shared strings appear in the saved workbook. -------
$ws.Range("C7").Value = "zona residencial"
$ws.Range("C8").Value = "zona educativa"
$ws.Range("C10").Value = "zona de turismo"
$ws.Range("C9").Value = "alto trafico"
$ws.Range("C15").Value = "alto nivel de comercio entre las 2 y 7pm"

# --- 4. Strip the D/E/F "documentation" columns from the display sheet,
#        keeping cell formatting (style) where the original had it. -------
$ws.Range("D1:D18").ClearContents()
$ws.Range("E1").ClearContents()
$ws.Range("E16").ClearContents()
$ws.Range("F1").ClearContents()
$ws.Range("E2:E15").ClearContents()
$ws.Range("E17:E18").ClearContents()
$ws.Range("F2:F18").ClearContents()

# Row 16 used to need extra height to fit the long "palabras clave" /
# "Refeferencias" text that has now been cleared; shrink it back down to
# the two wrapped lines that remain in the "descripcion" cell.
$ws.Rows.Item(16).RowHeight = 28.8

# Leave the cursor/selection on the display sheet, matching the saved file.
$ws.Select()
